# Trade #1 closed at 2026-02-17 19:43:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1299.99   # Current Capital
$wsSummary.Range("B4").Value = -0.01     # Total P&L $
$wsSummary.Range("B5").Value = -0.2      # Total P&L %
$wsSummary.Range("B6").Value = 1         # Total Trades
$wsSummary.Range("B8").Value = 1         # Losing Trades

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.98999999999999   # Capital
$wsStatus.Range("D4").Value = 1                   # Trades
$wsStatus.Range("E4").Value = -0.01                # P&L $
$wsStatus.Range("F4").Value = -0.01                # P&L %

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("G2").Value = 0.4                 # Exit Price
$wsAllTrades.Range("H2").Value = "CLOSED"            # Status
$wsAllTrades.Range("I2").Value = -2.439              # P&L %
$wsAllTrades.Range("J2").Value = -0.01               # P&L $
$wsAllTrades.Range("K2").Value = 99.98999999999999   # Capital After
$wsAllTrades.Range("P2").Value = "early_exit"        # Exit Reason
$wsAllTrades.Range("Q2").Value = 0.11                # Duration (min)

# --- MarketMaking sheet ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G2").Value = 0.4
$wsMM.Range("H2").Value = "CLOSED"
$wsMM.Range("I2").Value = -2.439
$wsMM.Range("J2").Value = -0.01
$wsMM.Range("K2").Value = 99.98999999999999
$wsMM.Range("P2").Value = "early_exit"
$wsMM.Range("Q2").Value = 0.11
